$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row|D|E|G  (empty field = cell unchanged)
$rows = @(
    "2|310.85|-1.67%|23",
    "3|48.57||23",
    "4|5.223|1.07%|23",
    "5|0.07772|-3.89%|23",
    "6|4.522|-0.02%|23",
    "7|1.292|18.53%|23",
    "8|1.561|-7.13%|23",
    "9|0.1237|-5.27%|23",
    "10|0.1919|-0.09%|23",
    "11|0.09282|-1.37%|23",
    "12|0.04545|6.39%|23",
    "13|0.1049|0.43%|23",
    "14|0.001287|-2.13%|23",
    "15|0.04193|-1.80%|23",
    "16|0.005834|-0.60%|23",
    "17|3.352|-1.39%|23",
    "18|2.386|-1.05%|23",
    "19||2.24%|23",
    "20|8.169|-1.78%|23",
    "21|0.1365|-1.49%|23",
    "22|0.3037|-3.47%|23",
    "23|0.001296|1.31%|23",
    "24|0.004217|-0.70%|23",
    "25|0.0001358|0.96%|23",
    "26|0.0003561|-95.19%|23",
    "27|||23",
    "28|||23",
    "29|||23",
    "30|||23",
    "31|||23",
    "32|||23",
    "33|||23",
    "34|||23",
    "35|||23",
    "36|||23",
    "37|||23",
    "38|0.02564|-5.31%|23",
    "39|0.05762|5.46%|23",
    "40|0.01037|92.62%|23",
    "41|0.008040|3.42%|23",
    "42|0.1420|-0.03%|23",
    "43|0.008407|14.02%|23",
    "44|0.008524|-0.73%|23",
    "45|0.3112|-0.86%|23",
    "46|0.00006914|1.72%|23",
    "47|0.00000000755|0.98%|23",
    "48|0.05567|-19.85%|23",
    "49|0.004025|0.98%|23",
    "50|0.00002113|0.98%|23",
    "51|0.0002013|0.98%|23"
)

foreach ($entry in $rows) {
    $parts = $entry.Split("|")
    $r = [int]$parts[0]
    $dVal = $parts[1]
    $eVal = $parts[2]
    $gVal = $parts[3]

    if ($dVal -ne "") {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
        $cell.ClearFormats()
    }

    if ($eVal -ne "") {
        $cell = $ws.Cells.Item($r, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $eVal
        $cell.ClearFormats()
    }

    if ($gVal -ne "") {
        $cell = $ws.Cells.Item($r, 7)
        $cell.NumberFormat = "@"
        $cell.Value = $gVal
        $cell.ClearFormats()
    }
}
